# fix bug convert speed to mm/s
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: drawing number / file name and cost value
$ws.Range("C2").Value = "MM00235120_Knotenblech.dxf"
$ws.Range("G2").Value = "2,82"

# Update row 3: drawing number / file name and cost value
$ws.Range("C3").Value = "L00040312_Halterung.dxf"
$ws.Range("G3").Value = "0,79"

# Remove the now-obsolete rows 4-8 entirely
$ws.Range("A4:G8").EntireRow.Delete()

# Widen column C to fit the new, longer file names
# (target stored width 33.42578125 chars; engine quantizes ColumnWidth to
# 1/6-character pixel steps, so feed it the value whose nearest snap is closest)
$ws.Columns.Item(3).ColumnWidth = 32.592447916666664

# Move the active selection to C4
$ws.Range("C4").Select()
